$d = $word.ActiveDocument

# Vertical tab character == Word's "line break" (<w:br/>) within a single run.
$brk = [char]11

# --- New paragraph 1: "I THINK THIS ONE IS CORRECT:" -----------------------
$last = $d.Paragraphs.Last
$r = $last.Range
$r.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.InsertAfter("I THINK THIS ONE IS CORRECT:")

# --- New paragraph 2: the new SQL query ------------------------------------
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$lines = @(
    "SELECT customer_id, first_name, last_name, total_cost FROM(",
    "SELECT customer_id, first_name, last_name, SUM(labor_hour), (labor_hour * labor_cost_per_hour) + part_cost as total_cost",
    "FROM simple_auto_shop.tbl_customer ",
    "NATURAL JOIN simple_auto_shop.tbl_customer_vehicle ",
    "NATURAL JOIN simple_auto_shop.tbl_vehicle_order ",
    "NATURAL JOIN simple_auto_shop.tbl_order_service ",
    "NATURAL JOIN simple_auto_shop.tbl_service",
    "NATURAL JOIN simple_auto_shop.tbl_rate",
    "WHERE tbl_service.rate_id = tbl_rate.rate_id",
    "ORDER BY SUM(labor_hour) DESC) a",
    ";"
)
$queryText = [string]::Join($brk, $lines)
$p2.Range.InsertAfter($queryText)

# --- New paragraph 3: empty -------------------------------------------------
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last

# --- New paragraph 4: "7)" ---------------------------------------------------
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$p4.Range.InsertAfter("7)")

# --- New paragraph 5: empty -------------------------------------------------
$p4.Range.InsertParagraphAfter()
